$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook already contains two "Rules String Hello"/"Hello2" decision
# tables (rows 3:11 and 27:35, each a 9-row block: title + header + 7 rows).
# This change adds a third such table ("Hello3") right below, at rows 39:47,
# by duplicating the first table's layout/styles/merged title cell and then
# renaming its title and re-attaching the same header/explanatory comments
# that live on the first table (B3,B4,C4,E4,C5,E5).

$src = $ws.Range("B3:E11")
$dst = $ws.Range("B39")
$src.Copy($dst)

# Give the new table its own title text (a new shared string).
$ws.Range("B39").Value2 = "Rules String Hello3 (Integer hour)"

# Re-create the explanatory comments on the new table's header cells,
# mirroring the ones found on the first table.
$ws.Range("B39").AddComment("This is so-called Decision Table Header. It starts with the keyword ""Rules"".") | Out-Null
$ws.Range("B40").AddComment("`nRule column header. Rule column is used to to name particular rule rows for documentation and tracing purposes. It is also useful to create rule rows that span more than one cell vertically (this will be explained in one of the next tutorials)`n") | Out-Null
$ws.Range("C40").AddComment("Condition column header. Must start with ""C""") | Out-Null
$ws.Range("E40").AddComment("Return column header. Must start with ""RET"".  ") | Out-Null
$ws.Range("C41").AddComment("Condition expression. Must have type boolean. As you can see condition uses parameter hour from Method Header and variable min that defines column data. When condition is evaluated for each row, the cell value from this row is assigned to variable min") | Out-Null
$ws.Range("E41").AddComment("This is return expression performed for the first row where all conditions have been satisfied. The variable greeting is substittuted with a cell value from the rule row") | Out-Null

# Leave the selection near the newly added table, like the author did.
$ws.Range("F33").Select()
